# Fixed repetition counters in generator
#
# The document renders witness-repetition counters as a subscript letter
# right after the "<witness><digit?>" locus reference, e.g. "1/5d9 2b)" or
# "1/7c6 b <<" - the trailing subscript letter used to be a plain Latin "b"
# but should be the Greek letter "β" (beta) instead. Replace every run whose
# text is exactly the subscript "b" with subscript "β", without disturbing
# any neighbouring run (in particular the subscript "2" digit run that
# sometimes immediately precedes it must stay a separate run).

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Font.Subscript = $true

$matches = @()
while ($find.Execute("b", $true, $false, $false, $false, $false, $true, 1, $true, "", 0)) {
    $r = $find.Parent
    $matches += , @($r.Start, $r.End)
}

foreach ($pair in $matches) {
    $start = $pair[0]
    $end = $pair[1]
    $target = $d.Range($start, $end)

    # If the character immediately before this run shares the exact same
    # formatting (e.g. another subscript run, such as the repetition
    # counter's leading "2"), a plain text swap would cause the engine to
    # silently coalesce the two runs into one. Nudge the neighbour's
    # formatting away momentarily so the runs stay distinct, then restore
    # it once our replacement run exists.
    $needsGuard = $false
    if ($start -gt 0) {
        $prevChar = $d.Range($start - 1, $start)
        if ($prevChar.Font.Subscript -eq $target.Font.Subscript) {
            $needsGuard = $true
        }
    }

    if ($needsGuard) {
        $prevChar = $d.Range($start - 1, $start)
        $prevChar.Bold = 1
        $target.Text = "β"
        $prevChar2 = $d.Range($start - 1, $start)
        $prevChar2.Bold = 0
    } else {
        $target.Text = "β"
    }
}
